# Apply weekly update to "Fruta, Vega Modelo de Temuco - Piña" sheet.
# A new record is inserted as row 744 (pushing the existing "Caramelo"
# quality-block rows 744-799 down to 745-800).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 744; this shifts rows 744:799 down to 745:800
# and grows the used range / dimension to A1:T800 automatically.
$ws.Rows.Item(744).Insert()

# Populate the new row 744 with the new record. Columns A,B,C,E,F,G,H,I,J,K,R
# carry the same constant values as every other row in this "Piña / Caramelo"
# block; columns D,L,M,N,O,P,Q,S,T hold the new weekly figures.
$ws.Cells.Item(744, 1).Value = 10
$ws.Cells.Item(744, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(744, 3).Value = "La Araucanía"
$ws.Cells.Item(744, 4).Value = 45166
$ws.Cells.Item(744, 5).Value = 9
$ws.Cells.Item(744, 6).Value = "Fruta"
$ws.Cells.Item(744, 7).Value = 100108
$ws.Cells.Item(744, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(744, 9).Value = 100108005
$ws.Cells.Item(744, 10).Value = "Piña"
$ws.Cells.Item(744, 11).Value = "Caramelo"
$ws.Cells.Item(744, 12).Value = "Primera"
$ws.Cells.Item(744, 13).Value = 65
$ws.Cells.Item(744, 14).Value = 25000
$ws.Cells.Item(744, 15).Value = 25000
$ws.Cells.Item(744, 16).Value = 25000
$ws.Cells.Item(744, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(744, 18).Value = "Ecuador"
$ws.Cells.Item(744, 19).Value = 2083
$ws.Cells.Item(744, 20).Value = 12
